$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '55.760.23'
Set-TextValue $ws.Range('E2') '  +2.77%  '
# Row 3
Set-TextValue $ws.Range('D3') '2.495.53'
Set-TextValue $ws.Range('E3') '  +7.67%  '
# Row 4
Set-TextValue $ws.Range('E4') '  +0.08%  '
# Row 5
Set-TextValue $ws.Range('D5') '480.45'
Set-TextValue $ws.Range('E5') '  +7.42%  '
# Row 6
Set-TextValue $ws.Range('D6') '138.98'
Set-TextValue $ws.Range('E6') '  +7.33%  '
# Row 7
Set-TextValue $ws.Range('E7') '  +0.49%  '
# Row 8
Set-TextValue $ws.Range('D8') '0.512'
Set-TextValue $ws.Range('E8') '  +7.51%  '
# Row 9
Set-TextValue $ws.Range('D9') '2.488.28'
Set-TextValue $ws.Range('E9') '  +10.09%  '
# Row 10
Set-TextValue $ws.Range('D10') '0.0985'
Set-TextValue $ws.Range('E10') '  +6.49%  '
# Row 11
Set-TextValue $ws.Range('D11') '5.45'
Set-TextValue $ws.Range('E11') '  +0.74%  '
# Row 12
Set-TextValue $ws.Range('D12') '0.326'
Set-TextValue $ws.Range('E12') '  +4.56%  '
# Row 13
Set-TextValue $ws.Range('E13') '  +0.35%  '
# Row 14
Set-TextValue $ws.Range('D14') '2.931.26'
Set-TextValue $ws.Range('E14') '  +7.85%  '
# Row 15
Set-TextValue $ws.Range('D15') '55.752.97'
Set-TextValue $ws.Range('E15') '  +2.76%  '
# Row 16
Set-TextValue $ws.Range('D16') '0.0000137'
Set-TextValue $ws.Range('E16') '  +13.03%  '
# Row 17
Set-TextValue $ws.Range('D17') '20.44'
Set-TextValue $ws.Range('E17') '  +8.27%  '
# Row 18
Set-TextValue $ws.Range('D18') '2.502.49'
Set-TextValue $ws.Range('E18') '  +7.98%  '
# Row 19
Set-TextValue $ws.Range('D19') '4.34'
Set-TextValue $ws.Range('E19') '  +6.71%  '
# Row 20
Set-TextValue $ws.Range('D20') '321.06'
Set-TextValue $ws.Range('E20') '  +6.91%  '
# Row 21
Set-TextValue $ws.Range('E21') '  +5.44%  '
# Row 22
Set-TextValue $ws.Range('E22') '  -0.15%  '
# Row 23
Set-TextValue $ws.Range('E23') '  +6.01%  '
# Row 24
Set-TextValue $ws.Range('D24') '57.87'
Set-TextValue $ws.Range('E24') '  +3.78%  '
# Row 25
Set-TextValue $ws.Range('B25') 'Binance-PegBSC-USD'
Set-TextValue $ws.Range('C25') 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Range('D25') '1.01'
Set-TextValue $ws.Range('E25') '  +0.39%  '
# Row 26
Set-TextValue $ws.Range('B26') 'Polygon'
Set-TextValue $ws.Range('C26') 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range('D26') '0.404'
Set-TextValue $ws.Range('E26') '  +8.83%  '
# Row 27
Set-TextValue $ws.Range('B27') 'Kaspa'
Set-TextValue $ws.Range('C27') 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D27') '0.163'
Set-TextValue $ws.Range('E27') '  +4.14%  '
# Row 28
Set-TextValue $ws.Range('E28') '  +8.21%  '
# Row 29
Set-TextValue $ws.Range('E29') '  +7.32%  '
# Row 30
Set-TextValue $ws.Range('D30') '0.0₃0770'
Set-TextValue $ws.Range('E30') '  +7.84%  '
# Row 31
Set-TextValue $ws.Range('E31') '  +0.45%  '
# Row 32
Set-TextValue $ws.Range('D32') '148.38'
Set-TextValue $ws.Range('E32') '  +1.14%  '
# Row 33
Set-TextValue $ws.Range('E33') '  +6.14%  '
# Row 34
Set-TextValue $ws.Range('E34') '  +8.94%  '
# Row 35
Set-TextValue $ws.Range('D35') '5.18'
Set-TextValue $ws.Range('E35') '  +10.13%  '
# Row 36
Set-TextValue $ws.Range('D36') '3.68'
Set-TextValue $ws.Range('E36') '  +1.72%  '
# Row 37
Set-TextValue $ws.Range('E37') '  +9.06%  '
# Row 38
Set-TextValue $ws.Range('D38') '0.843'
Set-TextValue $ws.Range('E38') '  +0.12%  '
# Row 39
Set-TextValue $ws.Range('D39') '34.27'
Set-TextValue $ws.Range('E39') '  +3.70%  '
# Row 40
Set-TextValue $ws.Range('D40') '0.998'
Set-TextValue $ws.Range('E40') '  +0.32%  '
# Row 41
Set-TextValue $ws.Range('D41') '0.611'
Set-TextValue $ws.Range('E41') '  +17.48%  '
# Row 42
Set-TextValue $ws.Range('E42') '  +10.30%  '
# Row 43
Set-TextValue $ws.Range('D43') '3.37'
Set-TextValue $ws.Range('E43') '  +6.71%  '
# Row 44
Set-TextValue $ws.Range('E44') '  +5.91%  '
# Row 45
Set-TextValue $ws.Range('D45') '10.15'
Set-TextValue $ws.Range('E45') '  -1.22%  '
# Row 46
Set-TextValue $ws.Range('D46') '1.971.36'
Set-TextValue $ws.Range('E46') '  +1.78%  '
# Row 47
Set-TextValue $ws.Range('E47') '  +9.13%  '
# Row 48
Set-TextValue $ws.Range('D48') '0.0222'
Set-TextValue $ws.Range('E48') '  +6.48%  '
# Row 49
Set-TextValue $ws.Range('D49') '249.21'
Set-TextValue $ws.Range('E49') '  +31.05%  '
# Row 50
Set-TextValue $ws.Range('D50') '4.47'
Set-TextValue $ws.Range('E50') '  +9.21%  '
# Row 51
Set-TextValue $ws.Range('D51') '17.44'
Set-TextValue $ws.Range('E51') '  +7.65%  '

Write-Host "Applied cryptos list update"
